# Update cryptos list - price (D) and volume(1h) (E) columns, plus
# swap rows 28/29 (Cosmos <-> Toncoin) and rows 32/33 (InjectiveProtocol <-> OKB)
#
# All of these columns hold plain text in the source workbook (prices such as
# "51.954.14" or "358.26" are strings, not numbers). Assigning a plain string
# via .Value lets Excel auto-detect/convert number-looking strings into real
# numbers, which would change the cell type. To keep them as text we force
# the cell's number format to Text ("@") before assigning, then restore the
# default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.947.25"
Set-TextValue $ws.Range("E2") "  +0.59%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.787.24"
Set-TextValue $ws.Range("E3") "  -0.78%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.07%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "358.07"
Set-TextValue $ws.Range("E5") "  +1.76%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "109.25"
Set-TextValue $ws.Range("E6") "  -2.56%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.564"
Set-TextValue $ws.Range("E7") "  +0.14%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  -0.05%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("E9") "  -0.74%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "39.97"
Set-TextValue $ws.Range("E10") "  -3.24%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("E11") "  +0.61%  "

# Row 12 - TRON
Set-TextValue $ws.Range("E12") "  +1.29%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "19.45"
Set-TextValue $ws.Range("E13") "  -1.74%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.57"
Set-TextValue $ws.Range("E14") "  -1.99%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.227.90"
Set-TextValue $ws.Range("E15") "  -0.98%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.812.36"
Set-TextValue $ws.Range("E16") "  -0.20%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.949"
Set-TextValue $ws.Range("E17") "  +7.66%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "51.886.19"
Set-TextValue $ws.Range("E18") "  +0.54%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "7.40"
Set-TextValue $ws.Range("E19") "  -0.87%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("E20") "  -1.34%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.97"
Set-TextValue $ws.Range("E21") "  -2.68%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("E22") "  -0.88%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "274.37"
Set-TextValue $ws.Range("E23") "  +1.75%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "70.25"
Set-TextValue $ws.Range("E24") "  +0.98%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("E25") "  -0.49%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("E26") "  +0.29%  "

# Row 27 - Dai
Set-TextValue $ws.Range("E27") "  +0.08%  "

# Row 28/29 swap: Cosmos <-> Toncoin
Set-TextValue $ws.Range("B28") "Toncoin"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D28") "2.27"
Set-TextValue $ws.Range("E28") "  +1.41%  "

Set-TextValue $ws.Range("B29") "Cosmos"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D29") "10.18"
Set-TextValue $ws.Range("E29") "  -1.24%  "

# Row 30 - Kaspa
Set-TextValue $ws.Range("E30") "  +4.16%  "

# Row 31 - VeChain
Set-TextValue $ws.Range("D31") "0.0465"
Set-TextValue $ws.Range("E31") "  +3.58%  "

# Row 32/33 swap: InjectiveProtocol <-> OKB
Set-TextValue $ws.Range("B32") "OKB"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D32") "51.55"
Set-TextValue $ws.Range("E32") "  +1.94%  "

Set-TextValue $ws.Range("B33") "InjectiveProtocol"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D33") "34.46"
Set-TextValue $ws.Range("E33") "  +1.97%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "5.71"
Set-TextValue $ws.Range("E34") "  -1.23%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0846"
Set-TextValue $ws.Range("E35") "  +3.23%  "

# Row 36 - RenderToken
Set-TextValue $ws.Range("D36") "5.30"
Set-TextValue $ws.Range("E36") "  +1.06%  "

# Row 37 - FirstDigitalUSD
Set-TextValue $ws.Range("E37") "  -0.03%  "

# Row 38 - LidoDAOToken
Set-TextValue $ws.Range("E38") "  +0.28%  "

# Row 39 - ARBITRUM
Set-TextValue $ws.Range("E39") "  -2.57%  "

# Row 40 - Celestia
Set-TextValue $ws.Range("D40") "18.02"
Set-TextValue $ws.Range("E40") "  -0.29%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("E41") "  +2.23%  "

# Row 42 - Stellar
Set-TextValue $ws.Range("E42") "  -1.27%  "

# Row 43 - WEMIXToken
Set-TextValue $ws.Range("E43") "  -1.75%  "

# Row 44 - Monero
Set-TextValue $ws.Range("D44") "122.19"
Set-TextValue $ws.Range("E44") "  -3.14%  "

# Row 45 - EnergySwap
Set-TextValue $ws.Range("D45") "22.00"
Set-TextValue $ws.Range("E45") "  -7.52%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.074.44"
Set-TextValue $ws.Range("E46") "  +0.01%  "

# Row 47 - NEARProtocol
Set-TextValue $ws.Range("E47") "  -1.86%  "

# Row 48 - ApeXProtocol
Set-TextValue $ws.Range("E48") "  -4.06%  "

# Row 49 - THORChain
Set-TextValue $ws.Range("E49") "  +1.55%  "

# Row 50 - SEI
Set-TextValue $ws.Range("D50") "0.930"
Set-TextValue $ws.Range("E50") "  -0.23%  "

# Row 51 - FraxShare
Set-TextValue $ws.Range("E51") "  +0.49%  "
